$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.01301432890758512
$ws.Range("C2").Value = 0.01320457232952306
$ws.Range("D2").Value = 0.01310945061855409
$ws.Range("E2").Value = 0.00009512171096897042

$ws.Range("B3").Value = 0.3235294117647059
$ws.Range("C3").Value = 0.3205741626794258
$ws.Range("D3").Value = 0.3220517872220658
$ws.Range("E3").Value = 0.001477624542640027

$ws.Range("B4").Value = 0.02502211550612915
$ws.Range("C4").Value = 0.02536437630134393
$ws.Range("D4").Value = 0.02519324590373654
$ws.Range("E4").Value = 0.00017113039760739
